$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "(주)디셈버앤컴퍼니"
$ws.Range("B7").Value = "핀트(fint) 프론트엔드(front-end) 개발자"
$ws.Range("C7").Value = "https://www.jobplanet.co.kr/job/search?posting_ids%5B%5D=1290918"
$ws.Range("D7").Value = "경력"
$ws.Range("E7").Value = "react,typescript,jotai,emotion,webpack,babel,vue,angular"

$ws.Range("A8").Value = "(주)버즈빌"
$ws.Range("B8").Value = "[광고 추천팀] 백엔드 개발자 (Engineering Manager, Team Lead)"
$ws.Range("C8").Value = "https://www.jobplanet.co.kr/job/search?posting_ids%5B%5D=1288726"
$ws.Range("D8").Value = "경력"
$ws.Range("E8").Value = "go,python,mysql,dynamodb,redis,elasticsearch,kafka,kubernetes,grpc,apache,prometheus,kubernetes,spinnaker,datadog,grafana,prometheus,loki,aws,gcp"
